# Cost Budget Estimate.xlsx - add Project Manager / Product Owner rows to the
# Labor Cost Estimate table, rebalance QA headcount, and drop the trailing
# blank row at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert two new rows right after "Developer" (row 7) and before the
#    "System Platform and Tools" section header (old row 8). Excel shifts all
#    formulas/merged-cells/references below automatically.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(9).Insert()

# Match the look of the existing Labor Cost Estimate rows (Quality Assurance /
# Developer, row 6) for the two new rows' boxed cells (C:F).
$ws.Range("C6:F6").Copy()
$ws.Range("C8:F8").PasteSpecial(-4122)
$ws.Range("C6:F6").Copy()
$ws.Range("C9:F9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Populate the new "Project Manager" / "Product Owner" rows.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value2 = "Project Manager"
$ws.Range("B8").Formula = "=B6"
$ws.Range("C8").Value2 = 100
$ws.Range("D8").Formula = "=B8*C8"
$ws.Range("E8").Value2 = 1
$ws.Range("F8").Formula = "=E8*D8"

$ws.Range("A9").Value2 = "Product Owner"
$ws.Range("B9").Formula = "=B8"
$ws.Range("C9").Value2 = 100
$ws.Range("D9").Formula = "=C9*B9"
$ws.Range("E9").Value2 = 1
$ws.Range("F9").Formula = "=E9*D9"

# Stray leftover total that ships with the edit (column I, row 9).
$ws.Range("I9").Value2 = 2342335.8199999998

# ---------------------------------------------------------------------------
# 3) The QA headcount (Count column) drops from 5 to 3; the Subtotal formula
#    in F6 recalculates automatically.
# ---------------------------------------------------------------------------
$ws.Range("E6").Value2 = 3

# ---------------------------------------------------------------------------
# 4) Update the Budget total so it covers the new rows too.
# ---------------------------------------------------------------------------
$ws.Range("B1").Formula = "=SUM(F6:F9,F12:F20)"

# ---------------------------------------------------------------------------
# 5) Drop the trailing blank row at the very bottom of the sheet (old row 25,
#    now shifted to row 27) and keep the lone formatted cell that used to sit
#    on old row 24 (now row 26) aligned with the row above it.
# ---------------------------------------------------------------------------
$ws.Range("G24").NumberFormat = $ws.Range("G26").NumberFormat()
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(26).Delete()

# ---------------------------------------------------------------------------
# 6) Misc view bookkeeping to match the saved workbook (cursor position).
# ---------------------------------------------------------------------------
$ws.Range("E26").Select()
